# Weekly price update: insert a new report row for
# "Comercializadora del Agro de Limarí" / Arveja Verde just above the
# existing row 64, pushing the following rows down by one and extending
# the sheet's used range from R92 to R93.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 64 (shifts rows 64:92 down to 65:93).
$ws.Rows(64).Insert()

# Populate the newly inserted row with this week's data.
$ws.Range("A64").Value = 2
$ws.Range("B64").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C64").Value = "Coquimbo"
$ws.Range("D64").Value = 45119
$ws.Range("E64").Value = 4
$ws.Range("F64").Value = 100112022
$ws.Range("G64").Value = "Arveja Verde"
$ws.Range("H64").Value = "Perfection"
$ws.Range("I64").Value = "Primera"
$ws.Range("J64").Value = 1000
$ws.Range("K64").Value = 26000
$ws.Range("L64").Value = 28000
$ws.Range("M64").Value = 27000
$ws.Range("N64").Value = "`$/malla 25 kilos"
$ws.Range("O64").Value = "Provincia de Limarí"
$ws.Range("P64").Value = 1080
$ws.Range("Q64").Value = 25
$ws.Range("R64").Value = "Hortaliza"
